$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row: "_old" suffix -> "_FV2310", "_new" suffix -> "_FV2404".
#    Column K ("diff") is left untouched.
$headerCols = @("A","B","C","D","E","F","G","H","I","J","L","M","N","O","P","Q","R","S","T","U")
foreach ($col in $headerCols) {
    $cell = $ws.Range($col + "1")
    $val = $cell.Value()
    if ($null -ne $val -and $val -like "*_old") {
        $cell.Value = ($val -replace "_old$", "_FV2310")
    } elseif ($null -ne $val -and $val -like "*_new") {
        $cell.Value = ($val -replace "_new$", "_FV2404")
    }
}

# 2. Freeze the header row (row 1) - select the first cell below the header
#    and freeze panes there.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the data range into an Excel Table ("Table1") with a header row + autofilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U75"), $null, 1)
$tbl.Name = "Table1"
